# NB14 erweitert + excel_output start und end_time2
#
# The "Kosten" sheet gets a second data row (time = 10s) and the existing
# row 2 values for "Beste Werte" / "Optimale Werte" are corrected. The
# line chart on the sheet needs its category/value series extended to
# cover the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kosten")

# --- update existing row 2 values ---
$ws.Range("B2").Value = 25350
$ws.Range("C2").Value = 14950

# --- add new row 3 ---
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = 14950
$ws.Range("C3").Value = 14950

# --- extend the chart series ranges to include the new row ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$ser1 = $chart.SeriesCollection().Item(1)
$ser1.Formula = "=SERIES(Kosten!`$B`$1,Kosten!`$A`$2:`$A`$3,Kosten!`$B`$2:`$B`$3,1)"

$ser2 = $chart.SeriesCollection().Item(2)
$ser2.Formula = "=SERIES(Kosten!`$C`$1,Kosten!`$A`$2:`$A`$3,Kosten!`$C`$2:`$C`$3,2)"
